# Daily attendance processing - 2025-10-13 18:52:06
# Rotate the "Recorded By" list in column G: move the first recorder in
# each comma-separated list to the end (oldest-first -> most-recent-last
# style re-ordering produced by the nightly attendance sync job).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Value2

    if ($text -and $text.ToString().Trim().Length -gt 0) {
        $rawParts = $text.ToString().Split(",")

        if ($rawParts.Count -gt 1) {
            $trimmed = @()
            foreach ($p in $rawParts) {
                $trimmed += $p.Trim()
            }

            $first = $trimmed[0]
            $rest = $trimmed[1..($trimmed.Count - 1)]
            $rotated = $rest + $first

            $newText = [string]::Join(", ", $rotated)
            $cell.Value = $newText
        }
    }
}
